$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6 (pushes old row 6 down to row 7)
$ws.Rows.Item(6).Insert()

# Update values per new ordering/values
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = -1.968

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = -1.383

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = -0.5669999999999999

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = -0.005

$ws.Range("A6").Value = 999
$ws.Range("B6").Value = 1.027

$ws.Range("A7").Value = 4
$ws.Range("B7").Value = 1.384
